$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($range, [string]$value)
    # Force the cell to remain text even when the string parses as a
    # number (e.g. "1.00", "315.48"), then drop the number-format
    # override so no stray style index is left behind on the cell.
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextCell "D2" "43.140.21"
$ws.Range("E2").Value = "  +0.54%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.583.01"
$ws.Range("E3").Value = "  +2.15%  "

# Row 4 - TetherUSD
Set-TextCell "D4" "1.00"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
Set-TextCell "D5" "315.48"
$ws.Range("E5").Value = "  -0.55%  "

# Row 6 - Solana
$ws.Range("E6").Value = "  +1.54%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.06%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +1.47%  "

# Row 10 - Avalanche
Set-TextCell "D10" "35.67"
$ws.Range("E10").Value = "  -1.05%  "

# Row 11 - Dogecoin
Set-TextCell "D11" "0.0814"
$ws.Range("E11").Value = "  +0.45%  "

# Row 12 - Polkadot
Set-TextCell "D12" "7.52"
$ws.Range("E12").Value = "  -0.93%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextCell "D13" "2.979.68"
$ws.Range("E13").Value = "  +2.10%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -3.22%  "

# Row 15 - WrappedEther
Set-TextCell "D15" "2.582.07"
$ws.Range("E15").Value = "  +1.61%  "

# Row 16 - Chainlink
$ws.Range("E16").Value = "  -0.74%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  -0.40%  "

# Row 18 - WrappedBTC
Set-TextCell "D18" "43.238.53"
$ws.Range("E18").Value = "  +0.55%  "

# Row 19 - Uniswap
Set-TextCell "D19" "6.84"
$ws.Range("E19").Value = "  +2.88%  "

# Row 20 - InternetComputer(DFINITY)
Set-TextCell "D20" "12.55"
$ws.Range("E20").Value = "  -3.33%  "

# Row 21 - ShibaInu
Set-TextCell "D21" "0.0₃0963"
$ws.Range("E21").Value = "  -0.26%  "

# Row 22 - Litecoin
Set-TextCell "D22" "69.50"
$ws.Range("E22").Value = "  -0.93%  "

# Row 23 - BitcoinCash
Set-TextCell "D23" "254.05"
$ws.Range("E23").Value = "  +0.97%  "

# Row 24 - PancakeSwap
Set-TextCell "D24" "2.97"
$ws.Range("E24").Value = "  +0.29%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  +2.98%  "

# Row 26 - EthereumClassic
Set-TextCell "D26" "27.34"
$ws.Range("E26").Value = "  +1.75%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  -0.02%  "

# Row 28 - Toncoin
$ws.Range("E28").Value = "  +1.71%  "

# Row 29 - InjectiveProtocol
Set-TextCell "D29" "40.18"
$ws.Range("E29").Value = "  +0.69%  "

# Row 30 - Cosmos
Set-TextCell "D30" "10.34"
$ws.Range("E30").Value = "  +0.56%  "

# Row 31 - Filecoin
$ws.Range("E31").Value = "  -2.84%  "

# Row 32 - Monero
Set-TextCell "D32" "155.32"
$ws.Range("E32").Value = "  +0.37%  "

# Row 33 - LidoDAOToken
$ws.Range("E33").Value = "  +3.99%  "

# Row 34 - ARBITRUM
$ws.Range("E34").Value = "  +2.19%  "

# Row 35 - Hedera
$ws.Range("E35").Value = "  +1.98%  "

# Row 36 - WEMIXToken
$ws.Range("E36").Value = "  +3.54%  "

# Row 37 - Celestia
Set-TextCell "D37" "18.71"
$ws.Range("E37").Value = "  -0.75%  "

# Row 39 - ApeXProtocol
Set-TextCell "D39" "2.49"
$ws.Range("E39").Value = "  +9.92%  "

# Row 40 - Stellar
$ws.Range("E40").Value = "  -0.48%  "

# Row 41 - EnergySwap
Set-TextCell "D41" "22.57"
$ws.Range("E41").Value = "  -4.55%  "

# Row 42 - RenderToken
$ws.Range("E42").Value = "  +4.40%  "

# Row 43 - VeChain
$ws.Range("E43").Value = "  +0.06%  "

# Row 44 - FirstDigitalUSD
$ws.Range("E44").Value = "  +0.05%  "

# Row 45 - NEARProtocol
Set-TextCell "D45" "3.25"
$ws.Range("E45").Value = "  -1.05%  "

# Row 46 - Maker
Set-TextCell "D46" "2.010.33"
$ws.Range("E46").Value = "  -0.47%  "

# Row 47 - FraxShare
Set-TextCell "D47" "8.97"
$ws.Range("E47").Value = "  +2.02%  "

# Row 48 - RocketPoolETH
Set-TextCell "D48" "2.833.25"
$ws.Range("E48").Value = "  +2.06%  "

# Row 49 - BitcoinSV
Set-TextCell "D49" "83.04"
$ws.Range("E49").Value = "  -3.23%  "

# Row 50 - ordi
Set-TextCell "D50" "75.15"
$ws.Range("E50").Value = "  +2.22%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  +2.05%  "
